$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values (32 rows x columns A:B) ---
$ws.Range("A1").Value = -0.2564814264887687
$ws.Range("B1").Value = 0.2558516398354769
$ws.Range("A2").Value = -0.1696260102948921
$ws.Range("B2").Value = 0.16808970845788185
$ws.Range("A3").Value = -0.1183716229345837
$ws.Range("B3").Value = 0.11791323408973398
$ws.Range("A4").Value = -0.10991323414310195
$ws.Range("B4").Value = 0.1095038171187035
$ws.Range("A5").Value = -0.10650381714911017
$ws.Range("B5").Value = 0.10511695966863943
$ws.Range("A6").Value = -0.005836374005506073
$ws.Range("B6").Value = 0.00578732239840285
$ws.Range("A7").Value = 0.004212677527374975
$ws.Range("B7").Value = -0.004215736164508321
$ws.Range("A8").Value = 0.014215736090643194
$ws.Range("B8").Value = -0.014227478771556168
$ws.Range("A9").Value = 0.01622747874102659
$ws.Range("B9").Value = -0.01624766713019632
$ws.Range("A10").Value = -0.020641735042667975
$ws.Range("B10").Value = 0.02063349895439437
$ws.Range("A11").Value = -0.017633498988407936
$ws.Range("B11").Value = 0.017621156495522428
$ws.Range("A12").Value = -0.014121156532603152
$ws.Range("B12").Value = 0.014038284548063107
$ws.Range("A13").Value = -0.010538284586536939
$ws.Range("B13").Value = 0.010506907041548885
$ws.Range("A14").Value = -0.002506907103784428
$ws.Range("B14").Value = 0.0025025469033845837
$ws.Range("A15").Value = -0.0015025469290428362
$ws.Range("B15").Value = 0.0015016907113523104
$ws.Range("A16").Value = 0.0004983092577846016
$ws.Range("B16").Value = -0.0004988362461295637
$ws.Range("A17").Value = 0.0024988362154996224
$ws.Range("B17").Value = -0.0025005020747723705
$ws.Range("A18").Value = -0.08022562957348
$ws.Range("B18").Value = 0.08003634210918165
$ws.Range("A19").Value = -0.07603634213248789
$ws.Range("B19").Value = 0.07459535059457334
$ws.Range("A20").Value = -0.00801727545198716
$ws.Range("B20").Value = 0.008005867336063233
$ws.Range("A21").Value = -0.004005867369459182
$ws.Range("B21").Value = 0.003999999966433521
$ws.Range("A22").Value = -0.04571808555257917
$ws.Range("B22").Value = 0.04550297304389961
$ws.Range("A23").Value = -0.040502973080034366
$ws.Range("B23").Value = 0.04009982705894011
$ws.Range("A24").Value = -0.02009982717634351
$ws.Range("B24").Value = 0.019999999881148867
$ws.Range("A25").Value = -0.09728058569587539
$ws.Range("B25").Value = 0.097154284921249
$ws.Range("A26").Value = -0.09465428495799522
$ws.Range("B26").Value = 0.09449156584770435
$ws.Range("A27").Value = -0.09199156588650093
$ws.Range("B27").Value = 0.09102764827145071
$ws.Range("A28").Value = -0.08902764831668986
$ws.Range("B28").Value = 0.08836768749769419
$ws.Range("A29").Value = -0.08136768757533197
$ws.Range("B29").Value = 0.08117634139900165
$ws.Range("A30").Value = -0.021176341757094708
$ws.Range("B30").Value = 0.021024130406112285
$ws.Range("A31").Value = -0.014024130489151077
$ws.Range("B31").Value = 0.014001294666787345
$ws.Range("A32").Value = -0.004001294765844321
$ws.Range("B32").Value = 0.0039999999322386515

# --- Column widths ---
# Target stored OOXML widths are A=15.7109375, B=16.42578125.
# The engine quantizes Columns.ColumnWidth to steps of 1/6 character
# (stored = round(ColumnWidth*6)/6 + 5/6), so these inputs land on the
# closest reachable stored widths (15.666666... and 16.5 respectively).
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Columns.Item(2).ColumnWidth = 15.67

Write-Output "done"
